# Regenerate save_data column "K" (G column) values.
# The workbook stores per-row stats where column G (header "K") previously
# held a "Strike#"-style count; it is being regenerated using the new K
# calculation (std/mean derived s_vals), so only column G values change.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$kValues = @{
    2  = 0
    3  = 2
    4  = 2
    5  = 0
    6  = 0
    7  = 1
    8  = 1
    9  = 0
    10 = 0
    11 = 1
    12 = 0
    13 = 1
    14 = 1
    15 = 1
    16 = 0
    17 = 0
    18 = 3
    19 = 1
    20 = 1
    21 = 1
    22 = 2
    23 = 0
    24 = 0
    25 = 0
    26 = 1
    27 = 1
    28 = 0
    29 = 1
    30 = 1
    31 = 1
    32 = 1
    33 = 1
    34 = 0
    35 = 0
    36 = 1
    37 = 4
    38 = 1
    39 = 3
    40 = 2
    41 = 1
    42 = 1
    43 = 0
    44 = 1
    45 = 0
    46 = 1
    47 = 1
    48 = 1
    49 = 1
    50 = 1
    51 = 1
    52 = 1
    53 = 1
    54 = 2
    55 = 2
}

foreach ($row in $kValues.Keys) {
    $ws.Range("G$row").Value = $kValues[$row]
}
